# The workbook contained an embedded picture ("Picture 1") anchored over
# the worksheet (xl/drawings/drawing1.xml, referenced from Sheet1 via
# <drawing r:id="rId1"/>). As part of optimizing the delete logic (and
# adding support for a two-condition "delete by condition" action), that
# picture/shape is removed from the worksheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}
